$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 45133
$ws.Range("K2").Value = 22000
$ws.Range("L2").Value = 22000
$ws.Range("M2").Value = 22000
$ws.Range("P2").Value = 1467

# Row 4
$ws.Range("D4").Value = 44749
$ws.Range("J4").Value = 90
$ws.Range("K4").Value = 17000
$ws.Range("L4").Value = 18000
$ws.Range("M4").Value = 17556
$ws.Range("P4").Value = 1170

# Row 6
$ws.Range("D6").Value = 45091
$ws.Range("J6").Value = 40
$ws.Range("K6").Value = 20000
$ws.Range("L6").Value = 22000
$ws.Range("M6").Value = 21000
$ws.Range("P6").Value = 1400

# Row 7
$ws.Range("D7").Value = 45119
$ws.Range("K7").Value = 20000
$ws.Range("L7").Value = 20000
$ws.Range("M7").Value = 20000
$ws.Range("P7").Value = 1333

# Row 9
$ws.Range("D9").Value = 45141
$ws.Range("J9").Value = 50
$ws.Range("K9").Value = 8500
$ws.Range("L9").Value = 9000
$ws.Range("M9").Value = 8800
$ws.Range("P9").Value = 587
